$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") mirroring the existing header style (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$hdr = $ws.Range("I1:J1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Data for new columns I (I0) and J (IF), rows 2-45
$data = @(
    @(7,7),
    @(8,8),
    @(7,7),
    @(8,8),
    @(7,7),
    @(8,8),
    @(9,9),
    @(7,8),
    @(6,7),
    @(3,4),
    @(9,9),
    @(5,6),
    @(7,7),
    @(1,1),
    @(4,5),
    @(8,9),
    @(10,10),
    @(8,8),
    @(8,8),
    @(5,5),
    @(8,9),
    @(6,7),
    @(8,9),
    @(9,9),
    @(8,8),
    @(8,9),
    @(7,9),
    @(7,7),
    @(4,6),
    @(9,9),
    @(9,9),
    @(8,8),
    @(6,8),
    @(5,7),
    @(3,6),
    @(8,8),
    @(5,7),
    @(6,9),
    @(1,5),
    @(1,3),
    @(1,3),
    @(1,3),
    @(1,2),
    @(1,2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
